$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2307692307692308
$ws.Range("C2").Value = 0.4848484848484849
$ws.Range("J2").Value = 0.02331002331002331
$ws.Range("P2").Value = 0.1888111888111888
$ws.Range("S2").Value = 0.07226107226107226
$ws.Range("C3").Value = 0.03225806451612903
$ws.Range("J3").Value = 0.04147465437788019
$ws.Range("P3").Value = 0.7327188940092166
$ws.Range("S3").Value = 0.1935483870967742
$ws.Range("J4").Value = 0.01612903225806452
$ws.Range("P4").Value = 0.6612903225806451
$ws.Range("S4").Value = 0.3225806451612903
$ws.Range("B6").Value = 0.06880733944954129
$ws.Range("D6").Value = 0.004587155963302753
$ws.Range("F6").Value = 0.01376146788990826
$ws.Range("J6").Value = 0.3348623853211009
$ws.Range("O6").Value = 0.01834862385321101
$ws.Range("Q6").Value = 0.1605504587155963
$ws.Range("R6").Value = 0.04587155963302753
$ws.Range("S6").Value = 0.3532110091743119
$ws.Range("B7").Value = 0.1203703703703704
$ws.Range("D7").Value = 0.009259259259259259
$ws.Range("F7").Value = 0.03703703703703703
$ws.Range("J7").Value = 0.1435185185185185
$ws.Range("O7").Value = 0.01851851851851852
$ws.Range("Q7").Value = 0.1712962962962963
$ws.Range("R7").Value = 0.09259259259259259
$ws.Range("S7").Value = 0.4074074074074074
$ws.Range("B8").Value = 0.1199226305609284
$ws.Range("D8").Value = 0.02321083172147002
$ws.Range("F8").Value = 0.03288201160541586
$ws.Range("J8").Value = 0.1411992263056093
$ws.Range("O8").Value = 0.02321083172147002
$ws.Range("Q8").Value = 0.160541586073501
$ws.Range("R8").Value = 0.09477756286266925
$ws.Range("S8").Value = 0.4042553191489361
$ws.Range("B9").Value = 0.1302083333333333
$ws.Range("D9").Value = 0.01041666666666667
$ws.Range("F9").Value = 0.0625
$ws.Range("J9").Value = 0.140625
$ws.Range("O9").Value = 0.04166666666666666
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.0625
$ws.Range("S9").Value = 0.3854166666666667
$ws.Range("B10").Value = 0.1239509360877986
$ws.Range("D10").Value = 0.02969657843770174
$ws.Range("E10").Value = 0.001291155584247902
$ws.Range("F10").Value = 0.05810200129115559
$ws.Range("J10").Value = 0.1291155584247902
$ws.Range("O10").Value = 0.02324080051646223
$ws.Range("Q10").Value = 0.2162685603615236
$ws.Range("R10").Value = 0.07876049063912201
$ws.Range("S10").Value = 0.3395739186571982
$ws.Range("G11").Value = 0.1240105540897098
$ws.Range("J11").Value = 0.1319261213720317
$ws.Range("K11").Value = 0.2137203166226913
$ws.Range("L11").Value = 0.503957783641161
$ws.Range("S11").Value = 0.02638522427440633
$ws.Range("G12").Value = 0.7323232323232324
$ws.Range("J12").Value = 0.2121212121212121
$ws.Range("K12").Value = 0.005050505050505051
$ws.Range("L12").Value = 0.01515151515151515
$ws.Range("S12").Value = 0.03535353535353535
$ws.Range("G13").Value = 0.5625
$ws.Range("J13").Value = 0.3541666666666667
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.03333333333333333
$ws.Range("H15").Value = 0.1291666666666667
$ws.Range("I15").Value = 0.05
$ws.Range("J15").Value = 0.3
$ws.Range("K15").Value = 0.07083333333333333
$ws.Range("M15").Value = 0.008333333333333333
$ws.Range("O15").Value = 0.04583333333333333
$ws.Range("S15").Value = 0.3625
$ws.Range("F16").Value = 0.02181818181818182
$ws.Range("H16").Value = 0.1418181818181818
$ws.Range("I16").Value = 0.07636363636363637
$ws.Range("J16").Value = 0.44
$ws.Range("K16").Value = 0.1018181818181818
$ws.Range("M16").Value = 0.003636363636363636
$ws.Range("N16").Value = 0.007272727272727273
$ws.Range("O16").Value = 0.05818181818181818
$ws.Range("S16").Value = 0.1490909090909091
$ws.Range("F17").Value = 0.02529182879377432
$ws.Range("H17").Value = 0.1614785992217899
$ws.Range("I17").Value = 0.07198443579766536
$ws.Range("J17").Value = 0.443579766536965
$ws.Range("K17").Value = 0.1031128404669261
$ws.Range("M17").Value = 0.01750972762645914
$ws.Range("O17").Value = 0.0622568093385214
$ws.Range("S17").Value = 0.1147859922178988
$ws.Range("F18").Value = 0.004761904761904762
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.06666666666666667
$ws.Range("J18").Value = 0.4095238095238095
$ws.Range("K18").Value = 0.1285714285714286
$ws.Range("M18").Value = 0.02380952380952381
$ws.Range("O18").Value = 0.08095238095238096
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.01665510062456627
$ws.Range("H19").Value = 0.2359472588480222
$ws.Range("I19").Value = 0.07564191533657183
$ws.Range("J19").Value = 0.3671061762664816
$ws.Range("K19").Value = 0.1138098542678695
$ws.Range("M19").Value = 0.02428868841082582
$ws.Range("O19").Value = 0.05135322692574601
$ws.Range("S19").Value = 0.1151977793199167
